$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Remember the "Outstanding" column (P) values before the shift, so they can
# be restored with clean literals afterwards.
$outstanding = @{}
for ($r = 3; $r -le 14; $r++) {
    $outstanding[$r] = $ws.Cells.Item($r, 16).Value()
}

# Insert a new (blank) column before column N (14th column) - shifts
# the old "In Advance" / Outstanding columns one place to the right.
$ws.Columns.Item(14).Insert()
$ws.Columns.Item(14).ColumnWidth = 9.14

# Re-assign the shifted "Outstanding" column (now Q) values explicitly so
# they keep their original, compact numeric literal formatting instead of
# whatever the raw column-shift produced.
for ($r = 3; $r -le 14; $r++) {
    $ws.Cells.Item($r, 17).Value = $outstanding[$r]
}

# Make "Repayment Schedule" the active / selected sheet with a new selection,
# instead of "Summary".
$ws.Activate()
$ws.Range("R6").Select()
